# Add support for ExcelIgnoreAttribute and throw for recursive structures
# This adds a new "IntValue" column (C) to the "Third Sheet" worksheet,
# with a header and four integer values underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Third Sheet")

$ws.Range("C1").Value = "IntValue"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4

$ws.Range("C1").Select()
